# tdf#141309 / tdf#142215 follow-up to tdf#140968:
# Add a new filtered data row (time value 0.5 / "c") to the autofilter
# range and extend the autofilter + filter criteria so the new row's
# value (0.500, once formatted with the column's "0.000" number format)
# is kept in the set of visible filter values, alongside the pre-existing
# 0.046 and 0.516 criteria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new row of data (A8: time value, B8: shared string "c") ---
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "c"

# --- Re-apply the autofilter over the now-larger A1:B8 range, restoring
#     the same three filter values (now formatted values, not raw
#     doubles, avoiding the double string<->number conversion bug) ---
$ws.AutoFilterMode = $false
$rng = $ws.Range("A1:B8")
$rng.AutoFilter(1, @("0.046", "0.500", "0.516"), 7)

# --- Keep the hidden _FilterDatabase defined name in sync with the
#     new autofilter range (must use the sheet-qualified name; the
#     unqualified short name resolves to an empty/unrelated Name object) ---
$wb.Names.Item("Munka1!_FilterDatabase").RefersTo = "=Munka1!`$A`$1:`$B`$8"

# --- Match the recorded UI selection after the edit ---
$ws.Range("C7").Select()
